$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 1274.4
$ws.Range("I12").Value = 1306.1428
$ws.Range("J12").Value = 1200.3334
$ws.Range("K12").Value = 1306.1428
$ws.Range("L12").Value = 1200.3334
$ws.Range("M12").Value = -1136.1428
$ws.Range("N12").Value = -1540.3334
# Row 111
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("M111").ClearContents()
# Row 116
$ws.Range("H116").Value = 15007.454
$ws.Range("I116").Value = 5111
$ws.Range("K116").Value = 5111
$ws.Range("M116").Value = -1669
# Row 137
$ws.Range("H137").Value = 3268.3076
$ws.Range("I137").Value = 2869
$ws.Range("J137").Value = 4022.5557
$ws.Range("K137").Value = 8607
$ws.Range("L137").Value = 12067.6671
$ws.Range("M137").Value = -6057
$ws.Range("N137").Value = -17167.6671
# Row 138
$ws.Range("H138").Value = 3209.6956
$ws.Range("I138").Value = 2149.1667
$ws.Range("J138").Value = 3584
$ws.Range("K138").Value = 6447.500100000001
$ws.Range("L138").Value = 10752
$ws.Range("M138").Value = -1307.500100000001
$ws.Range("N138").Value = -21032
# Row 141
$ws.Range("H141").Value = 1136.8
$ws.Range("I141").Value = 1136.8
$ws.Range("K141").Value = 3410.4
$ws.Range("M141").Value = 1769.6

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2030.5
$ws.Range("I2").Value = 2100.75
$ws.Range("J2").Value = 1749.5
$ws.Range("K2").Value = 2100.75
$ws.Range("L2").Value = 1749.5
$ws.Range("M2").Value = -1987.75
$ws.Range("N2").Value = -1975.5
# Row 45
$ws.Range("H45").Value = 2171.9167
$ws.Range("I45").Value = 2229.2222
$ws.Range("K45").Value = 2229.2222
$ws.Range("M45").Value = -1852.2222
# Row 52
$ws.Range("H52").Value = 72998
$ws.Range("J52").Value = 72998
$ws.Range("L52").Value = 72998
$ws.Range("N52").Value = -73634
# Row 60
$ws.Range("H60").Value = 25000
$ws.Range("I60").Value = 25000
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 25000
$ws.Range("L60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -24267
# Row 61
$ws.Range("H61").Value = 1933.3948
$ws.Range("I61").Value = 1680.3125
$ws.Range("J61").Value = 3283.1667
$ws.Range("K61").Value = 1680.3125
$ws.Range("L61").Value = 3283.1667
$ws.Range("M61").Value = -1468.3125
$ws.Range("N61").Value = -3707.1667
# Row 102
$ws.Range("H102").Value = 3968.9285
$ws.Range("I102").Value = 3415.2727
$ws.Range("K102").Value = 3415.2727
$ws.Range("M102").Value = -1793.2727
# Row 116
$ws.Range("H116").Value = 2030.5
$ws.Range("I116").Value = 2100.75
$ws.Range("J116").Value = 1749.5
$ws.Range("K116").Value = 2100.75
$ws.Range("L116").Value = 1749.5
$ws.Range("M116").Value = 193.25
$ws.Range("N116").Value = -6337.5
# Row 136
$ws.Range("H136").Value = 1933.3948
$ws.Range("I136").Value = 1680.3125
$ws.Range("J136").Value = 3283.1667
$ws.Range("K136").Value = 5040.9375
$ws.Range("L136").Value = 9849.500100000001
$ws.Range("M136").Value = -2490.9375
$ws.Range("N136").Value = -14949.5001

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2030.5
$ws.Range("I3").Value = 2100.75
$ws.Range("J3").Value = 1749.5
$ws.Range("K3").Value = 2100.75
$ws.Range("L3").Value = 1749.5
$ws.Range("M3").Value = -1986.75
$ws.Range("N3").Value = -1977.5
# Row 86
$ws.Range("H86").Value = 1529.7
$ws.Range("I86").Value = 1261.25
$ws.Range("K86").Value = 1261.25
$ws.Range("M86").Value = -138.25
# Row 89
$ws.Range("H89").Value = 1529.7
$ws.Range("I89").Value = 1261.25
$ws.Range("K89").Value = 6306.25
$ws.Range("M89").Value = -690.25
# Row 105
$ws.Range("H105").Value = 1101.96
$ws.Range("I105").Value = 1069.4546
$ws.Range("J105").Value = 1340.3334
$ws.Range("K105").Value = 1069.4546
$ws.Range("L105").Value = 1340.3334
$ws.Range("M105").Value = 677.5454
$ws.Range("N105").Value = -4834.3334

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 6
$ws.Range("H6").Value = 5000
$ws.Range("J6").Value = 5000
$ws.Range("L6").Value = 5000
$ws.Range("N6").Value = -5226
# Row 7
$ws.Range("H7").Value = 141.16667
$ws.Range("J7").Value = 256.66666
$ws.Range("L7").Value = 256.66666
$ws.Range("N7").Value = -482.66666
# Row 31
$ws.Range("H31").Value = 8893.682000000001
$ws.Range("I31").Value = 8528.385
$ws.Range("J31").Value = 9421.333000000001
$ws.Range("K31").Value = 8528.385
$ws.Range("L31").Value = 9421.333000000001
$ws.Range("M31").Value = -8233.385
$ws.Range("N31").Value = -10011.333
# Row 34
$ws.Range("H34").Value = 8893.682000000001
$ws.Range("I34").Value = 8528.385
$ws.Range("J34").Value = 9421.333000000001
$ws.Range("K34").Value = 8528.385
$ws.Range("L34").Value = 9421.333000000001
$ws.Range("M34").Value = -8326.385
$ws.Range("N34").Value = -9825.333000000001
# Row 35
$ws.Range("H35").Value = 22975
$ws.Range("I35").Value = 22975
$ws.Range("K35").Value = 22975
$ws.Range("M35").Value = -22681
# Row 41
$ws.Range("H41").Value = 5160
$ws.Range("I41").Value = 5160
$ws.Range("K41").Value = 5160
$ws.Range("M41").Value = -4732
# Row 86
$ws.Range("H86").Value = 404760.2
$ws.Range("J86").Value = 673333.7
$ws.Range("L86").Value = 673333.7
$ws.Range("N86").Value = -675579.7
# Row 89
$ws.Range("H89").Value = 404760.2
$ws.Range("J89").Value = 673333.7
$ws.Range("L89").Value = 3366668.5
$ws.Range("N89").Value = -3377900.5

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 38
$ws.Range("H38").Value = 10.777778
$ws.Range("I38").Value = 4.25
$ws.Range("J38").Value = 16
$ws.Range("K38").Value = 12.75
$ws.Range("L38").Value = 48
$ws.Range("M38").Value = 334.25
$ws.Range("N38").Value = -742

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 223.93939
$ws.Range("I2").Value = 164.68182
$ws.Range("J2").Value = 342.45456
$ws.Range("K2").Value = 164.68182
$ws.Range("L2").Value = 342.45456
$ws.Range("M2").Value = -51.68181999999999
$ws.Range("N2").Value = -568.45456
# Row 15
$ws.Range("H15").Value = 30564.25
$ws.Range("I15").Value = 24919
$ws.Range("J15").Value = 47500
$ws.Range("K15").Value = 24919
$ws.Range("L15").Value = 47500
$ws.Range("M15").Value = -24631
$ws.Range("N15").Value = -48076
# Row 75
$ws.Range("H75").Value = 56999
$ws.Range("J75").Value = 56999
$ws.Range("L75").Value = 56999
$ws.Range("N75").Value = -58747
# Row 78
$ws.Range("H78").Value = 56999
$ws.Range("J78").Value = 56999
$ws.Range("L78").Value = 170997
$ws.Range("N78").Value = -179733
# Row 81
$ws.Range("H81").Value = 30564.25
$ws.Range("I81").Value = 24919
$ws.Range("J81").Value = 47500
$ws.Range("K81").Value = 24919
$ws.Range("L81").Value = 47500
$ws.Range("M81").Value = -23921
$ws.Range("N81").Value = -49496
# Row 84
$ws.Range("H84").Value = 30564.25
$ws.Range("I84").Value = 24919
$ws.Range("J84").Value = 47500
$ws.Range("K84").Value = 74757
$ws.Range("L84").Value = 142500
$ws.Range("M84").Value = -69765
$ws.Range("N84").Value = -152484
# Row 92
$ws.Range("H92").Value = 51470.918
$ws.Range("J92").Value = 51470.918
$ws.Range("L92").Value = 51470.918
$ws.Range("N92").Value = -55214.918
# Row 102
$ws.Range("H102").Value = 1998.3914
$ws.Range("I102").Value = 1976.1333
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 1976.1333
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -354.1333
$ws.Range("N102").Value = -6244
# Row 122
$ws.Range("H122").Value = 2556.889
$ws.Range("I122").Value = 2251.5
$ws.Range("K122").Value = 6754.5
$ws.Range("M122").Value = -4304.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 80
$ws.Range("H80").Value = 60128
$ws.Range("J80").Value = 60128
$ws.Range("L80").Value = 60128
$ws.Range("N80").Value = -62374
# Row 83
$ws.Range("H83").Value = 60128
$ws.Range("J83").Value = 60128
$ws.Range("L83").Value = 180384
$ws.Range("N83").Value = -191616
# Row 134
$ws.Range("H134").Value = 39624.75
$ws.Range("J134").Value = 39624.75
$ws.Range("L134").Value = 39624.75
$ws.Range("N134").Value = -49764.75
# Row 136
$ws.Range("H136").Value = 2745.92
$ws.Range("I136").Value = 1560.579
$ws.Range("K136").Value = 4681.737
$ws.Range("M136").Value = -2131.737

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 61
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
# Row 122
$ws.Range("H122").Value = 1443.0869
$ws.Range("J122").Value = 1459.3334
$ws.Range("L122").Value = 4378.0002
$ws.Range("N122").Value = -9278.0002
